# Reassign the row data for the rows that were reshuffled in this update.
# (The underlying row/cell formatting does not change -- only the
# Artfynd record values that occupy each row are swapped around.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 gets the data that used to be in row 17
$ws.Range("A16").Value = 111815269
$ws.Range("B16").Value = 90666
$ws.Range("D16").Value = "LC"
$ws.Range("E16").Value = 4364
$ws.Range("F16").Value = "Dropptaggsvamp"
$ws.Range("G16").Value = "Hydnellum ferrugineum"
$ws.Range("H16").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("P16").Value = "åsele 1:1 (åsele 1:1), Ås lm"
$ws.Range("Q16").Value = 610053.7842541422
$ws.Range("R16").Value = 7121273.15248157
$ws.Range("S16").Value = 1
$ws.Range("Z16").Value = "18:27"
$ws.Range("AB16").Value = "18:27"

# Row 17 gets the data that used to be in row 22
$ws.Range("A17").Value = 111814478
$ws.Range("B17").Value = 77515
$ws.Range("D17").Value = "NT"
$ws.Range("E17").Value = 6425
$ws.Range("F17").Value = "Garnlav"
$ws.Range("G17").Value = "Alectoria sarmentosa"
$ws.Range("H17").Value = "(Ach.) Ach."
$ws.Range("P17").Value = "åsele 1:1 (åsele 1:1), Ås lm"
$ws.Range("Q17").Value = 610155.3487898401
$ws.Range("R17").Value = 7121461.207019502
$ws.Range("S17").Value = 1
$ws.Range("Z17").Value = "17:41"
$ws.Range("AB17").Value = "17:41"

# Row 19 gets the data that used to be in row 20
$ws.Range("A19").Value = 111814688
$ws.Range("B19").Value = 90087
$ws.Range("D19").Value = "LC"
$ws.Range("E19").Value = 3298
$ws.Range("F19").Value = "Trådticka"
$ws.Range("G19").Value = "Climacocystis borealis"
$ws.Range("H19").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("P19").Value = "åsele 1:1 (åsele 1:1), Ås lm"
$ws.Range("Q19").Value = 610011.2059644217
$ws.Range("R19").Value = 7121475.688616944
$ws.Range("S19").Value = 1
$ws.Range("Z19").Value = "17:55"
$ws.Range("AB19").Value = "17:55"

# Row 20 gets the data that used to be in row 21
$ws.Range("A20").Value = 111814591
$ws.Range("B20").Value = 77515
$ws.Range("D20").Value = "NT"
$ws.Range("E20").Value = 6425
$ws.Range("F20").Value = "Garnlav"
$ws.Range("G20").Value = "Alectoria sarmentosa"
$ws.Range("H20").Value = "(Ach.) Ach."
$ws.Range("P20").Value = "åsele 1:1 (åsele 1:1), Ås lm"
$ws.Range("Q20").Value = 610012.4812897337
$ws.Range("R20").Value = 7121464.398116477
$ws.Range("S20").Value = 1
$ws.Range("Z20").Value = "17:50"
$ws.Range("AB20").Value = "17:50"

# Row 21 gets the data that used to be in row 19
$ws.Range("A21").Value = 111815114
$ws.Range("B21").Value = 90660
$ws.Range("D21").Value = "NT"
$ws.Range("E21").Value = 4362
$ws.Range("F21").Value = "Blå taggsvamp"
$ws.Range("G21").Value = "Hydnellum caeruleum"
$ws.Range("H21").Value = "(Hornem.) P.Karst."
$ws.Range("P21").Value = "åsele 1:1, Ås lm"
$ws.Range("Q21").Value = 610384.0265214761
$ws.Range("R21").Value = 7121170.261031131
$ws.Range("S21").Value = 5
$ws.Range("Z21").Value = "18:19"
$ws.Range("AB21").Value = "18:19"

# Row 22 gets the data that used to be in row 16
$ws.Range("A22").Value = 111814104
$ws.Range("B22").Value = 56398
$ws.Range("D22").Value = "NT"
$ws.Range("E22").Value = 100109
$ws.Range("F22").Value = "Tretåig hackspett"
$ws.Range("G22").Value = "Picoides tridactylus"
$ws.Range("H22").Value = "(Linnaeus, 1758)"
$ws.Range("P22").Value = "åsele 1:1 (åsele 1:1), Ås lm"
$ws.Range("Q22").Value = 610154.5078508666
$ws.Range("R22").Value = 7121460.305022033
$ws.Range("S22").Value = 1
$ws.Range("Z22").Value = "17:23"
$ws.Range("AB22").Value = "17:23"

# Row 25 gets the data that used to be in row 28
$ws.Range("A25").Value = 112013696
$ws.Range("B25").Value = 86961
$ws.Range("D25").Value = "NT"
$ws.Range("E25").Value = 4962
$ws.Range("F25").Value = "Mjölsvärting"
$ws.Range("G25").Value = "Lyophyllum semitale"
$ws.Range("H25").Value = "(Fr. : Fr.) Kühner"
$ws.Range("P25").Value = "Spångmyran, Ås lm"
$ws.Range("Q25").Value = 610070.1349689787
$ws.Range("R25").Value = 7121402.360087069
$ws.Range("S25").Value = 25
$ws.Range("Z25").Value = "19:40"
$ws.Range("AB25").Value = "19:40"

# Row 26 gets the data that used to be in row 29
$ws.Range("A26").Value = 112013703
$ws.Range("B26").Value = 77515
$ws.Range("D26").Value = "NT"
$ws.Range("E26").Value = 6425
$ws.Range("F26").Value = "Garnlav"
$ws.Range("G26").Value = "Alectoria sarmentosa"
$ws.Range("H26").Value = "(Ach.) Ach."
$ws.Range("P26").Value = "Spångmyran, Ås lm"
$ws.Range("Q26").Value = 610144.4332068264
$ws.Range("R26").Value = 7121461.253672058
$ws.Range("S26").Value = 25
$ws.Range("Z26").Value = "19:28"
$ws.Range("AB26").Value = "19:28"

# Row 27 gets the data that used to be in row 25
$ws.Range("A27").Value = 112013697
$ws.Range("B27").Value = 89423
$ws.Range("D27").Value = "NT"
$ws.Range("E27").Value = 5432
$ws.Range("F27").Value = "Granticka"
$ws.Range("G27").Value = "Porodaedalea chrysoloma"
$ws.Range("H27").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("P27").Value = "Spångmyran, Ås lm"
$ws.Range("Q27").Value = 610102.0736959254
$ws.Range("R27").Value = 7121412.654772604
$ws.Range("S27").Value = 25
$ws.Range("Z27").Value = "19:35"
$ws.Range("AB27").Value = "19:35"

# Row 28 gets the data that used to be in row 32
$ws.Range("A28").Value = 112013704
$ws.Range("B28").Value = 81248
$ws.Range("D28").Value = "NT"
$ws.Range("E28").Value = 1312
$ws.Range("F28").Value = "Gammelgransskål"
$ws.Range("G28").Value = "Pseudographis pinicola"
$ws.Range("H28").Value = "(Nyl.) Rehm"
$ws.Range("P28").Value = "Spångmyran, Ås lm"
$ws.Range("Q28").Value = 610093.591720929
$ws.Range("R28").Value = 7121454.644715369
$ws.Range("S28").Value = 25
$ws.Range("Z28").Value = "19:49"
$ws.Range("AB28").Value = "19:49"

# Row 29 gets the data that used to be in row 33
$ws.Range("A29").Value = 112013699
$ws.Range("B29").Value = 77515
$ws.Range("D29").Value = "NT"
$ws.Range("E29").Value = 6425
$ws.Range("F29").Value = "Garnlav"
$ws.Range("G29").Value = "Alectoria sarmentosa"
$ws.Range("H29").Value = "(Ach.) Ach."
$ws.Range("P29").Value = "Spångmyran, Ås lm"
$ws.Range("Q29").Value = 610068.1736430819
$ws.Range("R29").Value = 7121408.394281525
$ws.Range("S29").Value = 25
$ws.Range("Z29").Value = "19:40"
$ws.Range("AB29").Value = "19:40"

# Row 30 gets the data that used to be in row 31
$ws.Range("A30").Value = 112013690
$ws.Range("B30").Value = 88489
$ws.Range("D30").Value = "NT"
$ws.Range("E30").Value = 1962
$ws.Range("F30").Value = "Vaddporing"
$ws.Range("G30").Value = "Anomoporia kamtschatica"
$ws.Range("H30").Value = "(Parmasto) Bondartseva"
$ws.Range("P30").Value = "Spångmyran, Ås lm"
$ws.Range("Q30").Value = 610051.8565798617
$ws.Range("R30").Value = 7121425.252971379
$ws.Range("S30").Value = 25
$ws.Range("Z30").Value = "19:43"
$ws.Range("AB30").Value = "19:43"

# Row 31 gets the data that used to be in row 27
$ws.Range("A31").Value = 112013700
$ws.Range("B31").Value = 77515
$ws.Range("D31").Value = "NT"
$ws.Range("E31").Value = 6425
$ws.Range("F31").Value = "Garnlav"
$ws.Range("G31").Value = "Alectoria sarmentosa"
$ws.Range("H31").Value = "(Ach.) Ach."
$ws.Range("P31").Value = "Spångmyran, Ås lm"
$ws.Range("Q31").Value = 610101.9650201321
$ws.Range("R31").Value = 7121415.702941997
$ws.Range("S31").Value = 25
$ws.Range("Z31").Value = "19:35"
$ws.Range("AB31").Value = "19:35"

# Row 32 gets the data that used to be in row 26
$ws.Range("A32").Value = 112013691
$ws.Range("B32").Value = 88489
$ws.Range("D32").Value = "NT"
$ws.Range("E32").Value = 1962
$ws.Range("F32").Value = "Vaddporing"
$ws.Range("G32").Value = "Anomoporia kamtschatica"
$ws.Range("H32").Value = "(Parmasto) Bondartseva"
$ws.Range("P32").Value = "Spångmyran, Ås lm"
$ws.Range("Q32").Value = 610134.4051595986
$ws.Range("R32").Value = 7121460.896015909
$ws.Range("S32").Value = 25
$ws.Range("Z32").Value = "19:29"
$ws.Range("AB32").Value = "19:29"

# Row 33 gets the data that used to be in row 30
$ws.Range("A33").Value = 112013698
$ws.Range("B33").Value = 77515
$ws.Range("D33").Value = "NT"
$ws.Range("E33").Value = 6425
$ws.Range("F33").Value = "Garnlav"
$ws.Range("G33").Value = "Alectoria sarmentosa"
$ws.Range("H33").Value = "(Ach.) Ach."
$ws.Range("P33").Value = "Spångmyran, Ås lm"
$ws.Range("Q33").Value = 610094.4326785516
$ws.Range("R33").Value = 7121455.546697079
$ws.Range("S33").Value = 25
$ws.Range("Z33").Value = "19:49"
$ws.Range("AB33").Value = "19:49"

